$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 1749
$ws.Range("E2").Value = 252
$ws.Range("F2").Value = 252
$ws.Range("G2").Value = 229
$ws.Range("H2").Value = 191
$ws.Range("I2").Value = 191
$ws.Range("K2").Value = 3868
$ws.Range("L2").Value = 2666
$ws.Range("M2").Value = 1202
$ws.Range("N2").Value = 1202
$ws.Range("P2").Value = 77
$ws.Range("Q2").Value = -471
$ws.Range("R2").Value = -922
$ws.Range("S2").Value = 1280
$ws.Range("T2").Value = 808
$ws.Range("U2").Value = -1279
$ws.Range("V2").Value = 2118
$ws.Range("W2").Value = 14.39
$ws.Range("X2").Value = 10.94
$ws.Range("AA2").Value = 221.91
$ws.Range("AB2").Value = 1466.58
$ws.Range("AC2").Value = 1776
$ws.Range("AE2").Value = 10796
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = 347
$ws.Range("AI2").Value = 20.86
$ws.Range("AJ2").Value = 10380000
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()
$ws.Range("Y2").ClearContents()
$ws.Range("Z2").ClearContents()
$ws.Range("AD2").ClearContents()
$ws.Range("AH2").ClearContents()

# Row 3
$ws.Range("D3").Value = 2313
$ws.Range("E3").Value = 269
$ws.Range("F3").Value = 269
$ws.Range("G3").Value = 262
$ws.Range("H3").Value = 202
$ws.Range("I3").Value = 203
$ws.Range("J3").Value = -1
$ws.Range("K3").Value = 7552
$ws.Range("L3").Value = 5010
$ws.Range("M3").Value = 2542
$ws.Range("N3").Value = 1786
$ws.Range("O3").Value = 756
$ws.Range("P3").Value = 108
$ws.Range("Q3").Value = -2261
$ws.Range("R3").Value = -387
$ws.Range("S3").Value = 3438
$ws.Range("T3").Value = 25
$ws.Range("U3").Value = -2286
$ws.Range("V3").Value = 4447
$ws.Range("W3").Value = 11.62
$ws.Range("X3").Value = 8.710000000000001
$ws.Range("Y3").Value = 13.58
$ws.Range("Z3").Value = 3.53
$ws.Range("AA3").Value = 197.06
$ws.Range("AB3").Value = 1543.18
$ws.Range("AC3").Value = 1438
$ws.Range("AD3").Value = 30.16
$ws.Range("AE3").Value = 11055
$ws.Range("AF3").Value = 3.92
$ws.Range("AG3").Value = 267
$ws.Range("AH3").Value = 0.62
$ws.Range("AI3").Value = 21.23
$ws.Range("AJ3").Value = 16155000

# Row 4
$ws.Range("D4").Value = 2764
$ws.Range("E4").Value = 404
$ws.Range("F4").Value = 404
$ws.Range("G4").Value = 345
$ws.Range("H4").Value = 263
$ws.Range("I4").Value = 272
$ws.Range("J4").Value = -9
$ws.Range("K4").Value = 10214
$ws.Range("L4").Value = 7359
$ws.Range("M4").Value = 2855
$ws.Range("N4").Value = 2108
$ws.Range("O4").Value = 747
$ws.Range("P4").Value = 108
$ws.Range("Q4").Value = -1433
$ws.Range("R4").Value = -870
$ws.Range("S4").Value = 1941
$ws.Range("T4").Value = 2
$ws.Range("U4").Value = -1435
$ws.Range("V4").Value = 4747
$ws.Range("W4").Value = 14.61
$ws.Range("X4").Value = 9.529999999999999
$ws.Range("Y4").Value = 13.99
$ws.Range("Z4").Value = 2.96
$ws.Range("AA4").Value = 257.71
$ws.Range("AB4").Value = 1806.57
$ws.Range("AC4").Value = 1686
$ws.Range("AD4").Value = 17.47
$ws.Range("AE4").Value = 13048
$ws.Range("AF4").Value = 2.26
$ws.Range("AG4").Value = 353
$ws.Range("AH4").Value = 1.2
$ws.Range("AI4").Value = 20.96
$ws.Range("AJ4").Value = 16155000

# Row 5
$ws.Range("D5").Value = 3308
$ws.Range("E5").Value = 237
$ws.Range("F5").Value = 237
$ws.Range("G5").Value = 888
$ws.Range("H5").Value = 650
$ws.Range("I5").Value = 543
$ws.Range("J5").Value = 107
$ws.Range("K5").Value = 9413
$ws.Range("L5").Value = 6827
$ws.Range("M5").Value = 2586
$ws.Range("N5").Value = 2586
$ws.Range("P5").Value = 162
$ws.Range("Q5").Value = -1575
$ws.Range("R5").Value = 1909
$ws.Range("S5").Value = -336
$ws.Range("T5").Value = 246
$ws.Range("U5").Value = -1820
$ws.Range("V5").Value = 5334
$ws.Range("W5").Value = 7.15
$ws.Range("X5").Value = 19.65
$ws.Range("Y5").Value = 23.14
$ws.Range("Z5").Value = 6.62
$ws.Range("AA5").Value = 264.01
$ws.Range("AB5").Value = 1479.65
$ws.Range("AC5").Value = 3362
$ws.Range("AD5").Value = 10.01
$ws.Range("AE5").Value = 16008
$ws.Range("AF5").Value = 2.1
$ws.Range("AG5").Value = 600
$ws.Range("AH5").Value = 1.78
$ws.Range("AI5").Value = 17.85
$ws.Range("AJ5").Value = 16155000
$ws.Range("O5").ClearContents()

# Row 6
$ws.Range("D6").Value = 5628
$ws.Range("E6").Value = 851
$ws.Range("F6").Value = 851
$ws.Range("G6").Value = 691
$ws.Range("H6").Value = 524
$ws.Range("I6").Value = 524
$ws.Range("K6").Value = 14901
$ws.Range("L6").Value = 11423
$ws.Range("M6").Value = 3478
$ws.Range("N6").Value = 3478
$ws.Range("P6").Value = 190
$ws.Range("Q6").Value = -1345
$ws.Range("R6").Value = -2664
$ws.Range("S6").Value = 4420
$ws.Range("T6").Value = 2094
$ws.Range("U6").Value = -3438
$ws.Range("V6").Value = 9104
$ws.Range("W6").Value = 15.12
$ws.Range("X6").Value = 9.31
$ws.Range("Y6").Value = 17.29
$ws.Range("Z6").Value = 4.31
$ws.Range("AA6").Value = 328.42
$ws.Range("AB6").Value = 1701.81
$ws.Range("AC6").Value = 3235
$ws.Range("AD6").Value = 8.27
$ws.Range("AE6").Value = 18316
$ws.Range("AF6").Value = 1.46
$ws.Range("AI6").Value = 21.74
$ws.Range("AJ6").Value = 18990164
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()

# Row 7
$ws.Range("D7").Value = 3949
$ws.Range("E7").Value = 669
$ws.Range("G7").Value = 396
$ws.Range("H7").Value = 307
$ws.Range("I7").Value = 306
$ws.Range("K7").Value = 15902
$ws.Range("L7").Value = 12251
$ws.Range("M7").Value = 3653
$ws.Range("N7").Value = 3658
$ws.Range("P7").Value = 190
$ws.Range("Q7").Value = 1400
$ws.Range("R7").Value = -133
$ws.Range("S7").Value = 610
$ws.Range("T7").Value = 643
$ws.Range("U7").Value = 986
$ws.Range("W7").Value = 16.93
$ws.Range("X7").Value = 7.78
$ws.Range("Y7").Value = 8.58
$ws.Range("Z7").Value = 2
$ws.Range("AA7").Value = 335.39
$ws.Range("AC7").Value = 1612
$ws.Range("AD7").Value = 16.87
$ws.Range("AE7").Value = 19261
$ws.Range("AF7").Value = 1.41
$ws.Range("AG7").Value = 525
$ws.Range("AH7").Value = 1.93
$ws.Range("AI7").Value = 32.56

# Row 8
$ws.Range("D8").Value = 10158
$ws.Range("E8").Value = 1544
$ws.Range("G8").Value = 1311
$ws.Range("H8").Value = 992
$ws.Range("I8").Value = 992
$ws.Range("K8").Value = 17151
$ws.Range("L8").Value = 12611
$ws.Range("M8").Value = 4540
$ws.Range("N8").Value = 4532
$ws.Range("P8").Value = 190
$ws.Range("Q8").Value = 1680
$ws.Range("R8").Value = -505
$ws.Range("S8").Value = -244
$ws.Range("T8").Value = 456
$ws.Range("U8").Value = 1273
$ws.Range("W8").Value = 15.2
$ws.Range("X8").Value = 9.76
$ws.Range("Y8").Value = 24.22
$ws.Range("Z8").Value = 6
$ws.Range("AA8").Value = 277.75
$ws.Range("AC8").Value = 5224
$ws.Range("AD8").Value = 5.21
$ws.Range("AE8").Value = 23867
$ws.Range("AF8").Value = 1.14
$ws.Range("AG8").Value = 725
$ws.Range("AH8").Value = 2.67
$ws.Range("AI8").Value = 13.88

# Row 9
$ws.Range("D9").Value = 8770
$ws.Range("E9").Value = 932
$ws.Range("G9").Value = 984
$ws.Range("H9").Value = 743
$ws.Range("I9").Value = 743
$ws.Range("K9").Value = 17760
$ws.Range("L9").Value = 12610
$ws.Range("M9").Value = 5151
$ws.Range("N9").Value = 5188
$ws.Range("P9").Value = 190
$ws.Range("Q9").Value = 1430
$ws.Range("R9").Value = -465
$ws.Range("S9").Value = -78
$ws.Range("T9").Value = 474
$ws.Range("U9").Value = 1644
$ws.Range("W9").Value = 10.63
$ws.Range("X9").Value = 8.470000000000001
$ws.Range("Y9").Value = 15.28
$ws.Range("Z9").Value = 4.25
$ws.Range("AA9").Value = 244.81
$ws.Range("AC9").Value = 3912
$ws.Range("AD9").Value = 6.95
$ws.Range("AE9").Value = 27320
$ws.Range("AF9").Value = 1
$ws.Range("AG9").Value = 725
$ws.Range("AH9").Value = 2.67
$ws.Range("AI9").Value = 18.54
